# Apply the "new version with timestamp" update:
#  - insert a new item row for AVEROCOXIB (alphabetically between AUGRAM and CETAL)
#  - insert a new item row for GARAMYCIN (alphabetically between GABALEPSY and INCONT)
#  - renumber the leading index column for every item row
#  - refresh the grand-total and the generated timestamp in the footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the AVEROCOXIB row right above the current CETAL row (row 9).
# ---------------------------------------------------------------------------
$ws.Rows(9).Insert()
$ws.Rows(9).RowHeight = 25.5

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A9").NumberFormat = "General"
$ws.Range("A9").Value = 3

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "AVEROCOXIB 90 MG 20 F.C. TABS."

$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "0:1"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"

$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "228.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "114.0000"

$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "0:1"

# ---------------------------------------------------------------------------
# 2) Insert the GARAMYCIN row right above the current INCONT row.
#    After step 1, GABALEPSY sits at row 16 and INCONT at row 17, so the new
#    row goes in at 17.
# ---------------------------------------------------------------------------
$ws.Rows(17).Insert()
$ws.Rows(17).RowHeight = 25.5

$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

$ws.Range("A17").NumberFormat = "General"
$ws.Range("A17").Value = 10

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "GARAMYCIN 0.1% OINT. 15 GM"

$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "4:0"

$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = "1"

$ws.Range("N17").NumberFormat = "@"
$ws.Range("N17").Value = "22.00"

$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "22.0000"

$ws.Range("Q17").NumberFormat = "@"
$ws.Range("Q17").Value = "1:0"

# ---------------------------------------------------------------------------
# 3) Renumber the index column (A7:A40) sequentially 1..34 for all item rows.
# ---------------------------------------------------------------------------
for ($r = 7; $r -le 40; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "General"
    $ws.Cells.Item($r, 1).Value = ($r - 6)
}

# ---------------------------------------------------------------------------
# 4) Refresh the grand-total (now at row 41) and the generated timestamp in
#    the footer (now at row 42).
# ---------------------------------------------------------------------------
$ws.Range("P41").NumberFormat = "General"
$ws.Range("P41").Value = 1402.1700000000001

$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "Tuesday, 29 July, 2025 12:52 PM"
